$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 with Luiz's data
$ws.Range("A2").Value = "Luiz"
$ws.Range("B2").Value = 500
$ws.Range("C2").Value = 500
$ws.Range("D2").Value = 500

# Add row 3 for Joao
$ws.Range("A3").Value = "Joao"
$ws.Range("B3").Value = 200
$ws.Range("C3").Value = 400
$ws.Range("D3").Value = 400

# Add row 4 for Pedro
$ws.Range("A4").Value = "Pedro"
$ws.Range("B4").Value = 900
$ws.Range("C4").Value = 200
$ws.Range("D4").Value = 100
